# edit.ps1 - applies the "kleine korrekturen vor der abgabe" changes.
#
# Summary of changes applied:
#  1) The three figure-caption "SEQ Abbildung" fields are converted from the
#     <w:fldSimple> shorthand serialization into the explicit
#     begin/instrText/separate/result/end <w:fldChar> run sequence.
#  2) The sentence about the active player's highlight animation is
#     reworded/reflowed into four runs.
#  3) The portrait screenshot (anchorId 708BD7D7) run gains <w:rPr><w:noProof/></w:rPr>.

$d = $word.ActiveDocument

function Convert-SeqFieldParagraph {
    param(
        [int]$ParaIndex,
        [string]$NewInnerXml
    )
    $p = $d.Paragraphs.Item($ParaIndex)
    $rng = $d.Range($p.Range.Start, $p.Range.End)
    $xml = "<w:p>$NewInnerXml</w:p>"
    $rng.InsertXML($xml)
}

# --- 1a) "Abbildung 1" caption (Dialog zur Konfiguration des Spiels) ---
$inner1 = '<w:r><w:t xml:space="preserve">Abbildung </w:t></w:r>' +
          '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
          '<w:r><w:instrText xml:space="preserve"> SEQ Abbildung \* ARABIC </w:instrText></w:r>' +
          '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
          '<w:r><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r>' +
          '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' +
          '<w:r><w:t xml:space="preserve">: Dialog zur Konfiguration </w:t></w:r>' +
          '<w:r w:rsidR="00984101"><w:t>de</w:t></w:r>' +
          '<w:r><w:t>s Spiels</w:t></w:r>'

# --- 1b) "Abbildung 2" caption (Standardansicht des Spielfeldes) ---
$inner2 = '<w:r><w:t xml:space="preserve">Abbildung </w:t></w:r>' +
          '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
          '<w:r><w:instrText xml:space="preserve"> SEQ Abbildung \* ARABIC </w:instrText></w:r>' +
          '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
          '<w:r><w:rPr><w:noProof/></w:rPr><w:t>2</w:t></w:r>' +
          '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' +
          '<w:r><w:t>: Standardansicht des Spielfeldes</w:t></w:r>'

# --- 1c) "Abbildung 3" caption (Ansicht des Spielfeldes auf einem Smartphone) ---
$inner3 = '<w:r><w:t xml:space="preserve">Abbildung </w:t></w:r>' +
          '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
          '<w:r><w:instrText xml:space="preserve"> SEQ Abbildung \* ARABIC </w:instrText></w:r>' +
          '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
          '<w:r><w:rPr><w:noProof/></w:rPr><w:t>3</w:t></w:r>' +
          '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' +
          '<w:r><w:t>: Ansicht des Spielfeldes auf einem Smartphone</w:t></w:r>'

# Locate the three caption paragraphs by their current text (robust against
# any paragraph numbering drift) and rewrite each one's run content.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -eq "Abbildung 1: Dialog zur Konfiguration des Spiels") {
        Convert-SeqFieldParagraph -ParaIndex $i -NewInnerXml $inner1
    }
    elseif ($t -eq "Abbildung 2: Standardansicht des Spielfeldes") {
        Convert-SeqFieldParagraph -ParaIndex $i -NewInnerXml $inner2
    }
    elseif ($t -eq "Abbildung 3: Ansicht des Spielfeldes auf einem Smartphone") {
        Convert-SeqFieldParagraph -ParaIndex $i -NewInnerXml $inner3
    }
}

# --- 2) Reflow the "gerade am Zug ist ..." sentence -----------------------
$findRng = $d.Content
$findRng.Find.Execute("gerade am Zug ist", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$spanStart = $findRng.Start

$endRng = $d.Content
$endRng.Find.Execute("r Zeit.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$spanEnd = $endRng.End

$targetRng = $d.Range($spanStart, $spanEnd)
$newSentenceXml = '<w:p>' +
    '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">gerade am Zug ist, wird mit der blauen Primärfarbe hervorgehoben und erhält </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">nach kurzer Zeit </w:t></w:r>' +
    '<w:r><w:t>eine pulsierende Animation</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
$targetRng.InsertXML($newSentenceXml)

# --- 3) Mark the portrait-screenshot drawing's run as NoProof -------------
# (adds <w:rPr><w:noProof/></w:rPr> to the <w:r> that hosts the
#  anchorId="708BD7D7" drawing, matching Word's behaviour when that picture
#  is touched again before saving.)
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $w = [int]([math]::Round($shp.Width * 12700))
    $h = [int]([math]::Round($shp.Height * 12700))
    if ($w -eq 1742257 -or ($w -gt 1742000 -and $w -lt 1742500 -and $h -gt 3497000 -and $h -lt 3498000)) {
        $shp.Range.NoProofing = $true
    }
}
